$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set Experimental boolean value "true" in B7
$ws.Range("B7").Value = "true"

# Update Date value in B8
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"
